$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-empty data row 6 with new measurement values
$ws.Range("D6").Value = 0.00114
$ws.Range("E6").Value = 0.00113
$ws.Range("F6").Value = 0.00141
$ws.Range("G6").Value = 0.00133
$ws.Range("H6").Value = "10000 hits"
$ws.Range("I6").Value = 0.00116
$ws.Range("J6").Value = 0.0012
$ws.Range("K6").Value = 0.00142
$ws.Range("L6").Value = 0.00131

# Change the unit labels from "[ms]" to "[s]" in row 3 and row 10 (D,E,F,G,I,J,K,L)
foreach ($r in 3,10) {
    foreach ($col in "D","E","F","G","I","J","K","L") {
        $ws.Range("$col$r").Value = "[s]"
    }
}

# Update the active cell selection in the bottom-right pane to Q11
$ws.Range("Q11").Select()
